# Applies numeric value corrections across multiple sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) as produced by the scheduled profit-sheet refresh.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 109
$ws.Range("I9").Value = 109
$ws.Range("K9").Value = 109
$ws.Range("M9").Value = 60
$ws.Range("H12").Value = 368.42856
$ws.Range("J12").Value = 107.5
$ws.Range("L12").Value = 107.5
$ws.Range("N12").Value = -447.5
$ws.Range("H19").Value = 1772.9166
$ws.Range("I19").Value = 889.5
$ws.Range("K19").Value = 889.5
$ws.Range("M19").Value = -714.5
$ws.Range("H41").Value = 578.6
$ws.Range("I41").Value = 267.2
$ws.Range("J41").Value = 890
$ws.Range("K41").Value = 267.2
$ws.Range("L41").Value = 890
$ws.Range("M41").Value = 172.8
$ws.Range("N41").Value = -1770
$ws.Range("H80").Value = 1962.375
$ws.Range("I80").Value = 1897.3334
$ws.Range("J80").Value = 2046
$ws.Range("K80").Value = 5692.0002
$ws.Range("L80").Value = 6138
$ws.Range("M80").Value = -4694.0002
$ws.Range("N80").Value = -8134
$ws.Range("H83").Value = 1962.375
$ws.Range("I83").Value = 1897.3334
$ws.Range("J83").Value = 2046
$ws.Range("K83").Value = 17076.0006
$ws.Range("L83").Value = 18414
$ws.Range("M83").Value = -12084.0006
$ws.Range("N83").Value = -28398
$ws.Range("H92").Value = 1026343.94
$ws.Range("I92").Value = 1368229.6
$ws.Range("K92").Value = 1368229.6
$ws.Range("M92").Value = -1366981.6
$ws.Range("H112").Value = 7049
$ws.Range("J112").Value = 7049
$ws.Range("L112").Value = 21147
$ws.Range("N112").Value = -23363
$ws.Range("H138").Value = 2805.5789
$ws.Range("J138").Value = 2737.8
$ws.Range("L138").Value = 8213.400000000001
$ws.Range("N138").Value = -18493.4

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7470
$ws.Range("I61").Value = 4400
$ws.Range("K61").Value = 4400
$ws.Range("M61").Value = -4188
$ws.Range("H88").Value = 3159.7856
$ws.Range("J88").Value = 4059.7144
$ws.Range("L88").Value = 4059.7144
$ws.Range("N88").Value = -4871.7144
$ws.Range("H91").Value = 3159.7856
$ws.Range("J91").Value = 4059.7144
$ws.Range("L91").Value = 4059.7144
$ws.Range("N91").Value = -6867.7144
$ws.Range("H136").Value = 7470
$ws.Range("I136").Value = 4400
$ws.Range("K136").Value = 13200
$ws.Range("M136").Value = -10650

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 96875.95
$ws.Range("I86").Value = 1616.5
$ws.Range("K86").Value = 1616.5
$ws.Range("M86").Value = -493.5
$ws.Range("H89").Value = 96875.95
$ws.Range("I89").Value = 1616.5
$ws.Range("K89").Value = 8082.5
$ws.Range("M89").Value = -2466.5
$ws.Range("H105").Value = 1889.0454
$ws.Range("I105").Value = 1889.0454
$ws.Range("K105").Value = 1889.0454
$ws.Range("M105").Value = -142.0454
$ws.Range("H134").Value = 7819.3237
$ws.Range("I134").Value = 7384.5864
$ws.Range("K134").Value = 22153.7592
$ws.Range("M134").Value = -19618.7592

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1500596.8
$ws.Range("I58").Value = 2718452.5
$ws.Range("J58").Value = 1697.2307
$ws.Range("K58").Value = 2718452.5
$ws.Range("L58").Value = 1697.2307
$ws.Range("M58").Value = -2718249.5
$ws.Range("N58").Value = -2103.2307
$ws.Range("H62").Value = 2950
$ws.Range("J62").Value = 2950
$ws.Range("L62").Value = 2950
$ws.Range("N62").Value = -4198
$ws.Range("H65").Value = 2950
$ws.Range("J65").Value = 2950
$ws.Range("L65").Value = 14750
$ws.Range("N65").Value = -20990
$ws.Range("H107").Value = 378.33334
$ws.Range("I107").Value = 378.33334
$ws.Range("K107").Value = 378.33334
$ws.Range("M107").Value = 1541.66666
$ws.Range("H136").Value = 1500596.8
$ws.Range("I136").Value = 2718452.5
$ws.Range("J136").Value = 1697.2307
$ws.Range("K136").Value = 8155357.5
$ws.Range("L136").Value = 5091.6921
$ws.Range("M136").Value = -8152807.5
$ws.Range("N136").Value = -10191.6921

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 200.83333
$ws.Range("I2").Value = 462.5
$ws.Range("J2").Value = 70
$ws.Range("K2").Value = 2775
$ws.Range("L2").Value = 420
$ws.Range("M2").Value = -2662
$ws.Range("N2").Value = -646
$ws.Range("H40").Value = 316.66666
$ws.Range("I40").Value = 150
$ws.Range("J40").Value = 400
$ws.Range("K40").Value = 600
$ws.Range("L40").Value = 1600
$ws.Range("M40").Value = -531
$ws.Range("N40").Value = -1738
$ws.Range("H107").Value = 1429.875
$ws.Range("J107").Value = 1429.875
$ws.Range("L107").Value = 4289.625
$ws.Range("N107").Value = -8129.625
$ws.Range("H130").Value = 1430
$ws.Range("I130").Value = 1430
$ws.Range("J130").Value = 0
$ws.Range("K130").Value = 4290
$ws.Range("L130").Value = 0
$ws.Range("M130").Value = 730
$ws.Range("N130").ClearContents()
$ws.Range("H131").Value = 770.1739
$ws.Range("J131").Value = 786.25287
$ws.Range("L131").Value = 2358.75861
$ws.Range("N131").Value = -12438.75861

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2525.5
$ws.Range("I102").Value = 2557.5386
$ws.Range("J102").Value = 2479.2222
$ws.Range("K102").Value = 2557.5386
$ws.Range("L102").Value = 2479.2222
$ws.Range("M102").Value = -935.5385999999999
$ws.Range("N102").Value = -5723.2222

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 8499.700000000001
$ws.Range("J40").Value = 14999.4
$ws.Range("L40").Value = 14999.4
$ws.Range("N40").Value = -15271.4
$ws.Range("H46").Value = 2399.7
$ws.Range("I46").Value = 700
$ws.Range("J46").Value = 2588.5557
$ws.Range("K46").Value = 700
$ws.Range("L46").Value = 2588.5557
$ws.Range("M46").Value = -512
$ws.Range("N46").Value = -2964.5557
$ws.Range("H132").Value = 1751.7273
$ws.Range("I132").Value = 1544.1177
$ws.Range("K132").Value = 4632.3531
$ws.Range("M132").Value = -2102.3531
$ws.Range("H134").Value = 51103.2
$ws.Range("J134").Value = 51103.2
$ws.Range("L134").Value = 51103.2
$ws.Range("N134").Value = -61243.2
$ws.Range("H136").Value = 2862.7368
$ws.Range("I136").Value = 1656.9286
$ws.Range("K136").Value = 4970.7858
$ws.Range("M136").Value = -2420.7858

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 49714
$ws.Range("J46").Value = 49714
$ws.Range("L46").Value = 49714
$ws.Range("N46").Value = -50176
$ws.Range("H113").Value = 1138.1
$ws.Range("I113").Value = 755.1667
$ws.Range("J113").Value = 1712.5
$ws.Range("K113").Value = 2265.5001
$ws.Range("L113").Value = 5137.5
$ws.Range("M113").Value = -95.5001000000002
$ws.Range("N113").Value = -9477.5
$ws.Range("H134").Value = 49714
$ws.Range("J134").Value = 49714
$ws.Range("L134").Value = 149142
$ws.Range("N134").Value = -154212
$ws.Range("H136").Value = 18520048
$ws.Range("I136").Value = 25253500
$ws.Range("J136").Value = 3056.125
$ws.Range("K136").Value = 75760500
$ws.Range("L136").Value = 9168.375
$ws.Range("M136").Value = -75757950
$ws.Range("N136").Value = -14268.375

Write-Host "Applied 184 cell changes."
